# Insert a new daily price record as row 79 (pushing the existing rows
# 79-179 down to 80-180) and populate it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("79:79").Insert()

$ws.Range("A79").Value = 3
$ws.Range("B79").Value = "Femacal de La Calera"
$ws.Range("C79").Value = "Coquimbo"
$ws.Range("D79").Value = 44413
$ws.Range("E79").Value = 5
$ws.Range("F79").Value = 100112032
$ws.Range("G79").Value = "Zapallo italiano"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 130
$ws.Range("K79").Value = 5500
$ws.Range("L79").Value = 6000
$ws.Range("M79").Value = 5808
$ws.Range("N79").Value = "$/caja 70 unidades"
$ws.Range("O79").Value = "Región de Arica y Parinacota"
$ws.Range("P79").Value = 83
$ws.Range("Q79").Value = 70
$ws.Range("R79").Value = "Hortaliza"
